# dodanie wykresu burndown, update taskow
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01_SprintBacklog")
$ws.Activate()

# --- widen column B (task status column gets longer "IN PROGRESS" values) ---
$ws.Range("B1").EntireColumn.ColumnWidth = 23.875

# --- Id column (A5:A20) gets sequential task numbers; Status (B) updates for a few rows ---
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 2
$ws.Range("A7").Value = 3
$ws.Range("A8").Value = 4
$ws.Range("A9").Value = 5
$ws.Range("A10").Value = 6
$ws.Range("A11").Value = 7
$ws.Range("A12").Value = 8
$ws.Range("A13").Value = 9
$ws.Range("A14").Value = 10
$ws.Range("A15").Value = 11
$ws.Range("A16").Value = 12
$ws.Range("A17").Value = 13
$ws.Range("A18").Value = 14
$ws.Range("A19").Value = 15
$ws.Range("A20").Value = 16

$ws.Range("B5").Value = "IN PROGRESS"
$ws.Range("B6").Value = "IN PROGRESS"
$ws.Range("B10").Value = "IN PROGRESS"

# --- hours remaining left on some tasks increased (D column - original size) ---
$ws.Range("D11").Value = 20
$ws.Range("D12").Value = 20
$ws.Range("D14").Value = 10
$ws.Range("D17").Value = 9

# --- "pozostalo [h]" (remaining hours, E column) filled in for every task ---
$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 2.5
$ws.Range("E8").Value = 2.5
$ws.Range("E9").Value = 2
$ws.Range("E10").Value = 0.5
$ws.Range("E11").Value = 20
$ws.Range("E12").Value = 20
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 10
$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 0.5
$ws.Range("E17").Value = 9
$ws.Range("E18").Value = 2
$ws.Range("E19").Value = 4
$ws.Range("E20").Value = 1

# --- totals row ---
$ws.Range("A21").Value = "SUMA"
$ws.Range("D21").Formula = "=SUM(D5:D20)"
$ws.Range("E21").Formula = "=SUM(E5:E20)"

# --- burndown data table (rows 29-65) ---
$ws.Range("A29").Value = "data"
$ws.Range("B29").Value = "pozostało [h]"

$ws.Range("A30").Value = 41273
$ws.Range("B30").Value = 80.5
$ws.Range("A31").Value = 41274
$ws.Range("B31").Value = 79

$startSerial = 41275
for ($i = 0; $i -lt 34; $i++) {
    $r = 32 + $i
    $ws.Cells.Item($r, 1).Value = $startSerial + $i
}

$ws.Range("A30:A65").NumberFormat = "mm-dd-yy"

# --- chart: burndown line chart ---
$co = $ws.ChartObjects().Add(1143000, 0, 2476500, 2600000)
$chart = $co.Chart
$chart.ChartType = "Line"
$chart.HasTitle = $true
$chart.ChartTitle.Text = "burndown"
$chart.SeriesCollection().NewSeries()
$ser = $chart.SeriesCollection().Item(1)
$ser.Name = "Sprint1"
$ser.XValues = $ws.Range("A30:A65")
$ser.Values = $ws.Range("B30:B65")
$chart.HasLegend = $true
$chart.Legend.Position = "Right"
$chart.Axes(1).CategoryType = "xlTimeScale"
$chart.Axes(2).MinimumScale = 0
$chart.Axes(2).MaximumScale = 81
